$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Korea" results column is being added to the report, to the left
# of the existing "AWS Korea" column (M). Inserting a whole column shifts
# the old M/N data into N/O and carries styles/formulas along.
$ws.Columns("M").Insert()

# Header for the new column.
$ws.Range("M5").Value = "Korea"

# New "In Seconds" raw measurements for the new Korea column.
$ws.Range("M7").Value = 19.1
$ws.Range("M8").Value = 14.8
$ws.Range("M10").Value = 14.3

# Mirror header into the "In Gbps" table below (row 16 copies row 5).
$ws.Range("M16").Formula = "=M5"

# New "In Gbps" computed results (100/seconds) for rows that now have data.
$ws.Range("M18").Formula = "=100/M7"
$ws.Range("M18").NumberFormat = "0.0"

$ws.Range("M19").Formula = "=100/M8"
$ws.Range("M19").NumberFormat = "0.0"

$ws.Range("M21").Formula = "=100/M10"
$ws.Range("M21").NumberFormat = "0.0"

# Restore the selection to match where the author ended up editing.
$ws.Range("M21").Select() | Out-Null
